$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-16 Wednesday" "2025-07-17 Thursday"

Replace-Text "44÷5=8, 4" "22÷6=3, 4"
Replace-Text "28÷7=4, 0" "84÷3=28, 0"
Replace-Text "36÷2=18, 0" "24÷2=12, 0"
Replace-Text "40÷8=5, 0" "42÷8=5, 2"
Replace-Text "70÷3=23, 1" "42÷6=7, 0"

Replace-Text "89÷7=12, 5" "58÷2=29, 0"
Replace-Text "29÷4=7, 1" "45÷5=9, 0"
Replace-Text "43÷6=7, 1" "47÷5=9, 2"
Replace-Text "84÷8=10, 4" "18÷8=2, 2"
Replace-Text "56÷9=6, 2" "23÷9=2, 5"

Replace-Text "49÷2=24, 1" "60÷8=7, 4"
Replace-Text "49÷8=6, 1" "42÷3=14, 0"
Replace-Text "68÷9=7, 5" "33÷2=16, 1"
Replace-Text "89÷2=44, 1" "87÷4=21, 3"
Replace-Text "54÷2=27, 0" "72÷2=36, 0"

Replace-Text "63÷3=21, 0" "75÷5=15, 0"
Replace-Text "42÷4=10, 2" "24÷2=12, 0"
Replace-Text "50÷7=7, 1" "22÷5=4, 2"
Replace-Text "27÷2=13, 1" "13÷4=3, 1"
Replace-Text "31÷4=7, 3" "98÷4=24, 2"

Replace-Text "20÷8=2, 4" "25÷5=5, 0"
Replace-Text "13÷7=1, 6" "73÷9=8, 1"
Replace-Text "70÷9=7, 7" "94÷4=23, 2"
Replace-Text "52÷8=6, 4" "51÷7=7, 2"
Replace-Text "67÷5=13, 2" "84÷2=42, 0"
